$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2.543602785316182, 0.7100761000916123, 0.07643685232798703, 0.01724538480356141, 0.002599618695145693, 5.057023726099231, 0.3476466353527314),
    @(2.47834464801889, 0.6476342069420298, 0.06961147884922525, 0.01685981551969906, 0.002608489123916296, 4.74064561398589, 0.3361113865778691),
    @(2.441371437843259, 0.6097883181083716, 0.06547316880389076, 0.01661991848199218, 0.002614202434003671, 4.546817989492382, 0.3292711945659192),
    @(2.427075240903548, 0.5944861663194274, 0.06379957101751188, 0.01652133745898432, 0.00261659807013047, 4.467921967524319, 0.3265440405375983),
    @(2.424747694393943, 0.5919524156736315, 0.06352243265820334, 0.01650491777385898, 0.002616999944796916, 4.45482628350635, 0.3260948193752284),
    @(2.44117552387371, 0.6095814655705567, 0.06545054679129692, 0.01661859234417218, 0.002614234468847161, 4.545753627789821, 0.329234172051315),
    @(2.520454560193741, 0.6884413495802733, 0.07407232037492406, 0.01711307878199797, 0.002602622031609891, 4.947836479994947, 0.3436185214222718),
    @(2.700853146996337, 0.8471843765046856, 0.09141553070045916, 0.01805901327652659, 0.002581952588578791, 5.74057702367395, 0.3737843825649207),
    @(2.849139058330479, 0.9665799135336215, 0.1044519694731179, 0.01874137260878417, 0.002568027801788198, 6.326838204855022, 0.3971920337543509),
    @(2.920142182224026, 1.021557372377231, 0.1104529609660005, 0.01904947214734065, 0.002561962395420862, 6.594653638503303, 0.40812268737929),
    @(2.94754952563153, 1.042476222468053, 0.1127360571919809, 0.0191658434105868, 0.002559703919978879, 6.696251890325186, 0.4123033783979793),
    @(2.941623574827076, 1.037966447751273, 0.1122438698896957, 0.01914079356030918, 0.002560188622113121, 6.674362439056324, 0.4114011321776445),
    @(2.92238652116356, 1.023276346064506, 0.1106405763389375, 0.01905905188996559, 0.002561775822409078, 6.603008411098926, 0.4084657986767866),
    @(2.910671290416417, 1.014291407390942, 0.1096599146480912, 0.01900894475883952, 0.002562753014562299, 6.559326388507372, 0.4066732519949028),
    @(2.844571055365122, 0.96300075731682, 0.1040612539313059, 0.01872119388954729, 0.00256842957691904, 6.309360201272426, 0.3964834470647105),
    @(2.804936162358842, 0.9317088743484874, 0.1006450975057049, 0.01854409962104242, 0.002571980642095286, 6.156316252101448, 0.3903052783619074),
    @(2.782472477481633, 0.9137728093590454, 0.09868683436683057, 0.01844202032520492, 0.002574048464898015, 6.068393770167575, 0.3867782851049668),
    @(2.774923599939484, 0.9077105099479468, 0.0980249223051004, 0.01840741937305079, 0.002574752956675595, 6.038642019758981, 0.3855886352597366),
    @(2.80912080268547, 0.9350334783485437, 0.1010080645601334, 0.0185629740766915, 0.002571600004495375, 6.172597091601119, 0.3909602033347284),
    @(2.928022716599457, 1.027588431400375, 0.1111112094678219, 0.01908306925908487, 0.002561308584807838, 6.623961699879374, 0.4093268440926465),
    @(3.008767812350129, 1.088664025373816, 0.1177765190782765, 0.01942125253089166, 0.002554805998771025, 6.920027596107502, 0.4215727692617861),
    @(2.965391486245835, 1.056011688627677, 0.1142132494775865, 0.01924090491936603, 0.002558256213801665, 6.761906408689754, 0.4150144225720851),
    @(2.807227921645108, 0.933530256857523, 0.100843949342547, 0.01855444176457333, 0.002571772008979528, 6.16523632466857, 0.3906640342063241),
    @(2.649325842989924, 0.8037740558602877, 0.0866741483459208, 0.01780546611076339, 0.002587321264718814, 5.525537408660426, 0.3654088970922942)
)

$cols = @("B", "C", "D", "E", "G", "I", "L")
$startRow = 2

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Range($cols[$j] + $row).Value = $data[$i][$j]
    }
}

$wb.Save()